# Apply the two text/formatting changes from the commit diff:
# 1. Table header cell "Characteristic" -> "Baseline Characteristics", and
#    the run is no longer bold.
# 2. Table cell "Nonwhite" -> "Non-white".

$d = $word.ActiveDocument

# --- Change 1: "Characteristic" -> "Baseline Characteristics" (un-bold) ---
$rng1 = $d.Content
$find1 = $rng1.Find
$find1.ClearFormatting()
$find1.Text = "Characteristic"
$find1.MatchWholeWord = $true
$find1.MatchCase = $true
$find1.Forward = $true
$find1.Wrap = 0
$found1 = $find1.Execute()
if ($found1) {
    $rng1.Text = "Baseline Characteristics"
    $rng1.Font.Bold = $false
}

# --- Change 2: "Nonwhite" -> "Non-white" ---
$rng2 = $d.Content
$find2 = $rng2.Find
$find2.ClearFormatting()
$find2.Text = "Nonwhite"
$find2.MatchWholeWord = $true
$find2.MatchCase = $true
$find2.Forward = $true
$find2.Wrap = 0
$found2 = $find2.Execute()
if ($found2) {
    $rng2.Text = "Non-white"
}
